$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("AC3").Value = 9.199999999999999
$ws.Range("AD3").Value = 15
$ws.Range("AG3").Value = 12
# Row 4
$ws.Range("F4").Value = 2
$ws.Range("I4").Value = 1000
$ws.Range("K4").Value = 1000
$ws.Range("R4").Value = 1.26
$ws.Range("S4").Value = 3.2
$ws.Range("W4").Value = 1.58
# Row 5
$ws.Range("Q5").Value = 2.36
# Row 6
$ws.Range("F6").Value = 1.78
$ws.Range("G6").Value = 1.84
$ws.Range("I6").Value = 5.4
$ws.Range("J6").Value = 3.75
$ws.Range("N6").Value = 3.85
$ws.Range("T6").Value = 1.81
$ws.Range("U6").Value = 2.06
$ws.Range("V6").Value = 1.22
$ws.Range("W6").Value = 2.18
# Row 7
$ws.Range("AH7").Value = 18.5
$ws.Range("H7").Value = 1.61
$ws.Range("I7").Value = 1.62
$ws.Range("N7").Value = 5.1
$ws.Range("T7").Value = 1.74
$ws.Range("V7").Value = 2.6
# Row 10
$ws.Range("F10").Value = 1.12
$ws.Range("K10").Value = 980
$ws.Range("V10").Value = 1.04
$ws.Range("W10").Value = 1.01
# Row 11
$ws.Range("AB11").Value = 1000
$ws.Range("AF11").Value = 1000
$ws.Range("AG11").Value = 1000
$ws.Range("AJ11").Value = 1000
$ws.Range("AK11").Value = 1000
$ws.Range("AL11").Value = 980
$ws.Range("W11").Value = 1.01
# Row 12
$ws.Range("AA12").Value = 980
$ws.Range("AE12").Value = 980
$ws.Range("AI12").Value = 980
$ws.Range("AL12").Value = 980
$ws.Range("K12").Value = 3.4
$ws.Range("P12").Value = 1.87
# Row 13
$ws.Range("AK13").Value = 980
$ws.Range("G13").Value = 2.06
$ws.Range("U13").Value = 1.84
$ws.Range("Z13").Value = 980
# Row 14
$ws.Range("AK14").Value = 1000
$ws.Range("AL14").Value = 980
$ws.Range("F14").Value = 1.47
$ws.Range("G14").Value = 1.53
# Row 15
$ws.Range("AJ15").Value = 980
$ws.Range("AK15").Value = 980
$ws.Range("F15").Value = 2.2
$ws.Range("K15").Value = 3.25
$ws.Range("N15").Value = 2.76
$ws.Range("Z15").Value = 980
# Row 16
$ws.Range("F16").Value = 2.62
$ws.Range("K16").Value = 3.15
# Row 17
$ws.Range("Q17").Value = 2.26
# Row 18
$ws.Range("AI18").Value = 48
$ws.Range("AJ18").Value = 140
$ws.Range("AK18").Value = 85
$ws.Range("AL18").Value = 100
$ws.Range("AM18").Value = 200
$ws.Range("AN18").Value = 130
$ws.Range("AO18").Value = 17.5
$ws.Range("L18").Value = 1.52
$ws.Range("N18").Value = 3.05
$ws.Range("O18").Value = 1.46
$ws.Range("P18").Value = 1.7
$ws.Range("Q18").Value = 2.36
$ws.Range("R18").Value = 1.25
$ws.Range("S18").Value = 4.7
$ws.Range("T18").Value = 2.16
$ws.Range("U18").Value = 1.81
$ws.Range("X18").Value = 10.5
$ws.Range("Y18").Value = 7.2
# Row 19
$ws.Range("AB19").Value = 1000
$ws.Range("AD19").Value = 1000
$ws.Range("AE19").Value = 980
$ws.Range("AF19").Value = 1000
$ws.Range("AH19").Value = 1000
$ws.Range("AJ19").Value = 980
$ws.Range("AK19").Value = 980
$ws.Range("AL19").Value = 980
$ws.Range("I19").Value = 3.8
$ws.Range("P19").Value = 1.89
$ws.Range("X19").Value = 1000
$ws.Range("Y19").Value = 1000
$ws.Range("Z19").Value = 980
# Row 22
$ws.Range("Q22").Value = 1.56
$ws.Range("R22").Value = 1.32
$ws.Range("S22").Value = 2.16
# Row 23
$ws.Range("L23").Value = 1.23
# Row 25
$ws.Range("AK25").Value = 980
$ws.Range("I25").Value = 3.45
$ws.Range("R25").Value = 1.23
# Row 26
$ws.Range("L26").Value = 1.22
# Row 27
$ws.Range("AB27").Value = 42
$ws.Range("AK27").Value = 140
$ws.Range("I27").Value = 1.41
$ws.Range("J27").Value = 5.2
$ws.Range("L27").Value = 1.26
$ws.Range("Q27").Value = 1.56
$ws.Range("V27").Value = 3.45
# Row 28
$ws.Range("AM28").Value = 130
$ws.Range("G28").Value = 1.84
$ws.Range("J28").Value = 3.8
$ws.Range("V28").Value = 1.24
$ws.Range("W28").Value = 2.18
# Row 29
$ws.Range("P29").Value = 1.69
# Row 30
$ws.Range("AA30").Value = 120
$ws.Range("Q30").Value = 1.86
# Row 32
$ws.Range("AF32").Value = 980
$ws.Range("AJ32").Value = 980
$ws.Range("AK32").Value = 980
$ws.Range("AL32").Value = 980
$ws.Range("AO32").Value = 980
$ws.Range("Z32").Value = 980
# Row 33
$ws.Range("AA33").Value = 980
$ws.Range("AE33").Value = 980
$ws.Range("H33").Value = 1.43
$ws.Range("Q33").Value = 1.61
$ws.Range("T33").Value = 1.64
$ws.Range("U33").Value = 1.68
$ws.Range("Z33").Value = 980
# Row 34
$ws.Range("H34").Value = 3.8
$ws.Range("I34").Value = 4.7
$ws.Range("Z34").Value = 40
# Row 35
$ws.Range("AC35").Value = 10
$ws.Range("AF35").Value = 980
$ws.Range("F35").Value = 2.36
$ws.Range("H35").Value = 2.82
$ws.Range("I35").Value = 3.7
$ws.Range("V35").Value = 1.44
$ws.Range("Z35").Value = 980
# Row 37
$ws.Range("AB37").Value = 10.5
$ws.Range("AC37").Value = 16.5
$ws.Range("I37").Value = 14.5
# Row 40
$ws.Range("I40").Value = 9.4
$ws.Range("J40").Value = 5.3
$ws.Range("K40").Value = 5.4
$ws.Range("O40").Value = 1.23
$ws.Range("T40").Value = 1.99
# Row 41
$ws.Range("S41").Value = 2.82
# Row 42
$ws.Range("AA42").Value = 980
$ws.Range("AE42").Value = 980
$ws.Range("Z42").Value = 980
# Row 43
$ws.Range("AK43").Value = 980
$ws.Range("AL43").Value = 980
$ws.Range("J43").Value = 3.4
$ws.Range("Z43").Value = 980
# Row 44
$ws.Range("J44").Value = 3.2
# Row 45
$ws.Range("F45").Value = 5
$ws.Range("G45").Value = 6
$ws.Range("H45").Value = 1.68
$ws.Range("I45").Value = 1.83
# Row 46
$ws.Range("F46").Value = 1.57
$ws.Range("G46").Value = 1.98
$ws.Range("H46").Value = 5
$ws.Range("K46").Value = 6.2
$ws.Range("L46").Value = 1.43
$ws.Range("N46").Value = 2.6
$ws.Range("P46").Value = 1.48
$ws.Range("Q46").Value = 2.16
$ws.Range("R46").Value = 1.12
$ws.Range("V46").Value = 1.11
$ws.Range("W46").Value = 2.02
# Row 47
$ws.Range("K47").Value = 4.1
# Row 48
$ws.Range("AL48").Value = 980
# Row 51
$ws.Range("F51").Value = 1.56
$ws.Range("G51").Value = 1.8
$ws.Range("H51").Value = 4.1
$ws.Range("I51").Value = 8.4
$ws.Range("J51").Value = 3.05
$ws.Range("W51").Value = 2.26
# Row 52
$ws.Range("G52").Value = 1.96
$ws.Range("K52").Value = 4.2
$ws.Range("W52").Value = 2.04
# Row 53
$ws.Range("L53").Value = 1.4
# Row 54
$ws.Range("F54").Value = 1.04
$ws.Range("H54").Value = 1.04
$ws.Range("J54").Value = 1.04
# Row 55
$ws.Range("AN55").Value = 230
$ws.Range("H55").Value = 1.49
# Row 56
$ws.Range("F56").Value = 2.32
